# Apply task-tracker updates to the "Tasks" worksheet:
#  - Row 11 (task #10): fill in the new "Add leagues" task
#      B11 = task text, C11 = LeagueController locations (multi-line),
#      D11 = Priority "Hight", E11 = Status "Opened"
#  - Row 8 (task #7): mark it "Closed" and add a Description note
#  - Move the active selection from B9 to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new task row -------------------------------------------------
# Copy cell formatting from existing rows that already carry the exact
# styles we need, so no brand-new style entries are introduced.
$ws.Range("B9").Copy()
$ws.Range("B11").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E6").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("B11").Value = "Реализовать добавление лиг на сайте. Редатирование, удаление."
$ws.Range("C11").Value = "LeagueController " + [char]10 + "Add" + [char]10 + "Edit" + [char]10 + "Delete"
$ws.Range("D11").Value = "Hight"
$ws.Range("E11").Value = "Opened"

$ws.Rows.Item(11).RowHeight = 60

# --- Row 8: close the task out and document the access rule --------------
$ws.Range("E7").Copy()
$ws.Range("E8").PasteSpecial(-4122)    # xlPasteFormats (Status -> "Closed" look)
$ws.Range("E8").Value = "Closed"
$ws.Range("F8").Value = "Статьи править и просматривать может только мастер. Админ может только просматривать."

# --- Selection update ------------------------------------------------------
$null = $ws.Range("B8").Select()
